# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '64.067.55'
Set-TextValue $ws.Range("E2") '  -1.15%  '
Set-TextValue $ws.Range("D3") '3.151.91'
Set-TextValue $ws.Range("E3") '  -0.98%  '
Set-TextValue $ws.Range("E4") '  -0.13%  '
Set-TextValue $ws.Range("D5") '601.52'
Set-TextValue $ws.Range("E5") '  -1.94%  '
Set-TextValue $ws.Range("D6") '142.22'
Set-TextValue $ws.Range("E6") '  -3.50%  '
Set-TextValue $ws.Range("E7") '  -0.07%  '
Set-TextValue $ws.Range("D8") '3.145.51'
Set-TextValue $ws.Range("E8") '  -1.01%  '
Set-TextValue $ws.Range("E9") '  -0.29%  '
Set-TextValue $ws.Range("E10") '  -2.93%  '
Set-TextValue $ws.Range("D11") '5.40'
Set-TextValue $ws.Range("E11") '  -2.05%  '
Set-TextValue $ws.Range("E12") '  -2.21%  '
Set-TextValue $ws.Range("D13") '0.0000255'
Set-TextValue $ws.Range("E13") '  -2.69%  '
Set-TextValue $ws.Range("D14") '34.92'
Set-TextValue $ws.Range("E14") '  -3.57%  '
Set-TextValue $ws.Range("D15") '3.665.87'
Set-TextValue $ws.Range("E15") '  -1.13%  '
Set-TextValue $ws.Range("E16") '  +2.69%  '
Set-TextValue $ws.Range("D17") '64.015.94'
Set-TextValue $ws.Range("E17") '  -1.40%  '
Set-TextValue $ws.Range("D18") '3.139.14'
Set-TextValue $ws.Range("E18") '  -1.39%  '
Set-TextValue $ws.Range("E19") '  -1.67%  '
Set-TextValue $ws.Range("D20") '487.56'
Set-TextValue $ws.Range("E20") '  +0.82%  '
Set-TextValue $ws.Range("E21") '  -0.71%  '
Set-TextValue $ws.Range("E22") '  -1.85%  '
Set-TextValue $ws.Range("E23") '  -2.86%  '
Set-TextValue $ws.Range("D24") '88.45'
Set-TextValue $ws.Range("E24") '  +4.64%  '
Set-TextValue $ws.Range("D25") '13.22'
Set-TextValue $ws.Range("E25") '  -4.96%  '
Set-TextValue $ws.Range("E26") '  +0.07%  '
Set-TextValue $ws.Range("E27") '  -2.64%  '
Set-TextValue $ws.Range("D28") '8.22'
Set-TextValue $ws.Range("E28") '  -6.99%  '
Set-TextValue $ws.Range("D29") '6.98'
Set-TextValue $ws.Range("E29") '  -2.73%  '
Set-TextValue $ws.Range("D30") '2.07'
Set-TextValue $ws.Range("E30") '  -3.34%  '
Set-TextValue $ws.Range("D31") '27.51'
Set-TextValue $ws.Range("E31") '  +2.68%  '
Set-TextValue $ws.Range("E32") '  -6.47%  '
Set-TextValue $ws.Range("E33") '  -0.12%  '
Set-TextValue $ws.Range("D34") '2.65'
Set-TextValue $ws.Range("E34") '  -3.61%  '
Set-TextValue $ws.Range("E35") '  -3.16%  '
Set-TextValue $ws.Range("E36") '  +0.19%  '
Set-TextValue $ws.Range("D37") '52.76'
Set-TextValue $ws.Range("E38") '  -6.64%  '
Set-TextValue $ws.Range("E39") '  -8.85%  '
Set-TextValue $ws.Range("E40") '  -1.26%  '
Set-TextValue $ws.Range("D41") '432.32'
Set-TextValue $ws.Range("E41") '  -7.56%  '
Set-TextValue $ws.Range("E42") '  -0.68%  '
Set-TextValue $ws.Range("E43") '  -0.40%  '
Set-TextValue $ws.Range("D44") '2.911.89'
Set-TextValue $ws.Range("E44") '  +1.44%  '
Set-TextValue $ws.Range("E45") '  -4.26%  '
Set-TextValue $ws.Range("E46") '  -7.00%  '
Set-TextValue $ws.Range("E47") '  -2.40%  '
Set-TextValue $ws.Range("D49") '25.81'
Set-TextValue $ws.Range("E49") '  -4.33%  '
Set-TextValue $ws.Range("E50") '  +0.00%  '
Set-TextValue $ws.Range("D51") '120.88'
Set-TextValue $ws.Range("E51") '  -0.18%  '
